# Insert a new data row at row 187 (pushing the existing rows 187:250 down
# to 188:251) and populate the new row with its own Betarraga/Maule record.
# This mirrors the target diff: dimension grows from A1:R250 to A1:R251 and
# every row from 187 downward holds the data that used to belong to the row
# above it, with the brand-new row 187 carrying fresh values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 187:250 down to 188:251, leaving a blank row 187 behind.
$ws.Rows.Item(187).Insert()

# Populate the new row 187.
$ws.Cells.Item(187, 1).Value2  = 5
$ws.Cells.Item(187, 2).Value2  = 'Macroferia Regional de Talca'
$ws.Cells.Item(187, 3).Value2  = 'Maule'
$ws.Cells.Item(187, 4).Value2  = 44588
$ws.Cells.Item(187, 5).Value2  = 7
$ws.Cells.Item(187, 6).Value2  = 100114014
$ws.Cells.Item(187, 7).Value2  = 'Betarraga'
$ws.Cells.Item(187, 8).Value2  = 'Sin especificar'
$ws.Cells.Item(187, 9).Value2  = 'Primera'
$ws.Cells.Item(187, 10).Value2 = 3000
$ws.Cells.Item(187, 11).Value2 = 700
$ws.Cells.Item(187, 12).Value2 = 700
$ws.Cells.Item(187, 13).Value2 = 700
$ws.Cells.Item(187, 14).Value2 = '$/paquete 5 unidades'
$ws.Cells.Item(187, 15).Value2 = 'Región del Maule'
$ws.Cells.Item(187, 16).Value2 = 140
$ws.Cells.Item(187, 17).Value2 = 5
$ws.Cells.Item(187, 18).Value2 = 'Hortaliza'
